$wb = $excel.ActiveWorkbook

# --- Rename the two placeholder sheets to their real test-data names -----
$wsShoeSize = $wb.Worksheets.Item("ESPN")
$wsShoeSize.Name = "shoeSize"

$wsBrand = $wb.Worksheets.Item("Verizon")
$wsBrand.Name = "Brand"

# --- Populate the shoeSize sheet with a header + numeric shoe sizes ------
$wsShoeSize.Range("A1").Value = "shoeSize"
$wsShoeSize.Range("A2").Value = 6
$wsShoeSize.Range("A2").NumberFormat = "0"
$wsShoeSize.Range("A3").Value = 6.5
$wsShoeSize.Range("A4").Value = 7
$wsShoeSize.Range("A5").Value = 7.5
$wsShoeSize.Range("A6").Value = 8
$wsShoeSize.Range("A7").Value = 8.5
$wsShoeSize.Range("A8").Value = 9
$wsShoeSize.Range("A9").Value = 10

# --- Populate the Brand sheet with a header + brand names -----------------
$wsBrand.Range("A1").Value = "Brand"
$wsBrand.Range("A2").Value = "adidas"
$wsBrand.Range("A3").Value = "Clarks"
$wsBrand.Range("A4").Value = "Crocs"
$wsBrand.Range("A5").Value = "Michael Kors"
$wsBrand.Range("A6").Value = "Nike"
$wsBrand.Range("A7").Value = "Skechers"
$wsBrand.Range("A8").Value = "Unbranded"
$wsBrand.Range("A9").Value = "VANS"

# --- Page setup for the shoeSize sheet (portrait) --------------------------
$wsShoeSize.PageSetup.Orientation = 1

# --- Recreate the end-user selection state ---------------------------------
# Condtion was the previously active/selected tab; after the edits the
# shoeSize sheet has cell A2 selected, and Brand ends up as the active tab
# with the whole of column A selected.
$wsShoeSize.Activate() | Out-Null
$wsShoeSize.Range("A2").Select() | Out-Null

$wsBrand.Activate() | Out-Null
$wsBrand.Columns("A:A").Select() | Out-Null

Write-Host "Applied ebay test-case edits"
